$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.109.00'
$ws.Range("E2").Value = '  -3.26%  '
$ws.Range("D3").Value = '2.655.12'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.66%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.37%  '
$ws.Range("E10").Value = '  -3.81%  '
$ws.Range("E11").Value = '  -2.37%  '
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '3.123.08'
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").Value = '59.114.38'
$ws.Range("E14").Value = '  -3.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000137'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").Value = '2.646.78'
$ws.Range("E17").Value = '  -7.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '340.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("E19").Value = '  -4.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.414'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("D27").Value = '0.0₃0805'
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '148.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.903'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.882'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("E39").Value = '  -5.77%  '
$ws.Range("E40").Value = '  -3.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.619'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '275.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0973'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0535'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").Value = '2.030.80'
$ws.Range("E48").Value = '  -5.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0229'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.95%  '

Write-Host "Applied crypto updates"